$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.512.82'
$ws.Range("E2").Value = '  -1.38%  '

# Row 3
$ws.Range("D3").Value = '1.748.79'
$ws.Range("E3").Value = '  -1.47%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.61'
$ws.Range("E5").Value = '  +0.40%  '

# Row 6
$ws.Range("E6").Value = '  -0.18%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4481'
$ws.Range("E7").Value = '  +4.89%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3600'
$ws.Range("E8").Value = '  -0.51%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07500'
$ws.Range("E9").Value = '  +0.20%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.95'
$ws.Range("E10").Value = '  -5.81%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.094'
$ws.Range("E11").Value = '  -1.30%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.43%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.70'
$ws.Range("E13").Value = '  -4.14%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.023'
$ws.Range("E14").Value = '  -1.85%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.123'
$ws.Range("E15").Value = '  -2.65%  '

# Row 16
$ws.Range("D16").Value = '1.753.15'
$ws.Range("E16").Value = '  -2.08%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.54'
$ws.Range("E17").Value = '  +2.59%  '

# Row 18
$ws.Range("E18").Value = '  -0.04%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06387'
$ws.Range("E19").Value = '  +0.70%  '

# Row 20
$ws.Range("E20").Value = '  +0.25%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.80'
$ws.Range("E21").Value = '  -2.58%  '

# Row 22
$ws.Range("E22").Value = '  -1.51%  '

# Row 23
$ws.Range("D23").Value = '27.563.50'
$ws.Range("E23").Value = '  -1.24%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.24'
$ws.Range("E24").Value = '  -1.23%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.092'
$ws.Range("E25").Value = '  -2.17%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.32'
$ws.Range("E26").Value = '  +1.60%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.46'
$ws.Range("E27").Value = '  +0.71%  '

# Row 28
$ws.Range("D28").Value = '1.951.04'
$ws.Range("E28").Value = '  -1.99%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.086'
$ws.Range("E29").Value = '  -3.98%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.70'
$ws.Range("E30").Value = '  -0.26%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.084'
$ws.Range("E31").Value = '  -6.98%  '

# Row 32
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09087'
$ws.Range("E32").Value = '  +1.63%  '

# Row 33
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.665'
$ws.Range("E33").Value = '  +4.81%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.546'
$ws.Range("E34").Value = '  -2.44%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.98'
$ws.Range("E35").Value = '  -4.89%  '

# Row 36
$ws.Range("E36").Value = '  -1.32%  '

# Row 37
$ws.Range("B37").Value = 'TheSandbox'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6366'
$ws.Range("E37").Value = '  -1.31%  '

# Row 38
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06016'
$ws.Range("E38").Value = '  -0.57%  '

# Row 39
$ws.Range("E39").Value = '  -0.99%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.940'
$ws.Range("E40").Value = '  -2.36%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.203'
$ws.Range("E41").Value = '  +1.54%  '

# Row 42
$ws.Range("E42").Value = '  -1.16%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.766'
$ws.Range("E43").Value = '  -1.53%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.21'
$ws.Range("E44").Value = '  -3.78%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.721'
$ws.Range("E45").Value = '  +0.32%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5897'
$ws.Range("E46").Value = '  -1.27%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.38'
$ws.Range("E47").Value = '  -1.51%  '

# Row 48
$ws.Range("E48").Value = '  -1.43%  '

# Row 49
$ws.Range("E49").Value = '  +0.28%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06853'
$ws.Range("E50").Value = '  -0.81%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.28'
$ws.Range("E51").Value = '  -3.01%  '
